# "Added preliminary teams to group list"
#
# - renames Sheet1/2/3 to Group List / Java Team / Web Team
# - adds Strength / Weakness / Lead? columns (G/H/I) to the Group List sheet
# - populates the (previously empty) Java Team and Web Team sheets

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

$ws1.Name = "Group List"
$ws2.Name = "Java Team"
$ws3.Name = "Web Team"

# ---------------------------------------------------------------------------
# Group List sheet: new Stregnth / Weakeness / Lead? columns
# ---------------------------------------------------------------------------

$ws1.Range("G1").Value = "Stregnth"
$ws1.Range("H1").Value = "Weakeness"
$ws1.Range("I1").Value = "Lead?"
$ws1.Range("G1:I1").Font.Size = 18

# Ben Dudley
$ws1.Range("G2").Value = "Web/DB/Basic Java"
$ws1.Range("H2").Value = "Java"
$ws1.Range("I2").Value = "DB (y) /Web"

# David Fairbrother
$ws1.Range("G3").Value = "Java "
$ws1.Range("H3").Value = "Web"
$ws1.Range("I3").Value = "-"

# Jonathan Englund
$ws1.Range("G4").Value = "Java (Rusty) / Web"
$ws1.Range("H4").Value = "-"
$ws1.Range("I4").Value = "-"

# Joshua Doyle
$ws1.Range("G5").Value = "Java"
$ws1.Range("H5").Value = "Web"
$ws1.Range("I5").Value = "Java (y)"

# Liam Fitzgerald
$ws1.Range("G6").Value = "Web/DB?"
$ws1.Range("H6").Value = "Java"
$ws1.Range("I6").Value = "-"

# Maurice Corriette
$ws1.Range("G7").Value = "Testing/?"
$ws1.Range("H7").Value = "?"
$ws1.Range("I7").Value = "?"

# Oliver Earl
$ws1.Range("G8").Value = "Web/Small Java"
$ws1.Range("H8").Value = "-"
$ws1.Range("I8").Value = "?"

# Tim Anderson
$ws1.Range("G9").Value = "Web?"
$ws1.Range("H9").Value = "Java?"
$ws1.Range("I9").Value = "?"

$ws1.Range("G2:I9").Font.Size = 12

# column widths
$ws1.Columns.Item(3).ColumnWidth = 18.15
$ws1.Columns.Item(4).ColumnWidth = 23.3
$ws1.Columns.Item(7).ColumnWidth = 20.5

$ws1.Range("G6").Select() | Out-Null

# ---------------------------------------------------------------------------
# Java Team sheet
# ---------------------------------------------------------------------------

$ws2.Range("A1").Value = "Joshua Doyle"
$ws2.Range("B1").Value = "Team Leader"
$ws2.Range("A2").Value = "David Fairbrother"
$ws2.Range("B2").Value = "Team Members"
$ws2.Range("A3").Value = "Jonathan Englund"
$ws2.Range("B3").Value = "Team Members"
$ws2.Range("A4").Value = "Ben Dudley?"
$ws2.Range("B4").Value = "Team Members?"

$ws2.Columns.Item(1).ColumnWidth = 14.33
$ws2.Columns.Item(2).ColumnWidth = 13.83

$ws2.Range("B5").Select() | Out-Null

# ---------------------------------------------------------------------------
# Web Team sheet
# ---------------------------------------------------------------------------

$ws3.Range("A1").Value = "Ben Dudley?"
$ws3.Range("B1").Value = "Team Leader?"
$ws3.Range("A2").Value = "Liam Fitzgerald"
$ws3.Range("B2").Value = "Team member"
$ws3.Range("A3").Value = "Oliver Earl "
$ws3.Range("B3").Value = "Team member?"
$ws3.Range("A4").Value = "Tim Anderson"
$ws3.Range("B4").Value = "Team member?"

$ws3.Columns.Item(1).ColumnWidth = 12.33
$ws3.Columns.Item(2).ColumnWidth = 13.0

$ws3.Range("E8").Select() | Out-Null

$ws1.Activate() | Out-Null
